# ---------------------------------------------------------------------------
# completezza funzionale_frontend.xlsx  -  apply the "Bug fix / Gestione del
# token / Pagina avvisi migliorata UI / Some changes signed in todotree.txt"
# commit.
#
# Strategy: the shared-string table is rebuilt by the engine on save by
# scanning cells and appending any still-referenced legacy strings (in their
# original relative order) followed by brand-new strings in the order they
# are first *written* by this script. We therefore issue our writes in the
# exact order needed to reproduce the target shared-string ordering.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header block (rows 2-4): "DISEGNO XML PAGINE STATICHE" section becomes
#    the "Dispensa" section, with two new sub-steps added under it.
# ---------------------------------------------------------------------------
$ws.Range("A13").ClearContents()                 # drop "RESPONSABILE:" row label

$ws.Range("B19").Value = "Login "                # new string idx 25
$ws.Range("B20").Value = "SignUp "                # new string idx 26

$ws.Range("D18").Value = "Visualizzazione piatti (VP)"          # idx 27
$ws.Range("E18").Value = "Struttura FXML statico"               # idx 28
$ws.Range("B18").Value = "Struttura FXML statico delle pagine"  # idx 29
$ws.Range("C18").Value = "Struttura FXML statico della pagina"  # idx 30

$ws.Range("D20").Value = "VP : aggiunta card visuale dinamica"                              # idx 31
$ws.Range("D21").Value = "VP : aggiunta funzionalità bottoni di cancellazione piatto"        # idx 32
$ws.Range("D22").Value = "VP : aggiunta funzionalità ricerca piatto"                         # idx 33
$ws.Range("D23").Value = "Ordina menu (OM)"                                                  # idx 34
$ws.Range("D25").Value = "OM : Aggiunta funzionalità scelta ordinamento"                     # idx 35
$ws.Range("D27").Value = "OM : Aggiunta funzionalità visualizza menu già presente"           # idx 36
$ws.Range("D26").Value = "OM : Aggiunta funzionalità reset categoria"                        # idx 37
$ws.Range("D28").Value = "Modifica Piatto (MP)"                                              # idx 38
$ws.Range("D30").Value = "MP: Inserimento dati e richieste"                                  # idx 39

$ws.Range("B4").Value  = "Aggiunta card visuale dinamica"                # idx 40
$ws.Range("E19").Value = "Aggiunta funzionalità scrittura avviso"        # idx 41
$ws.Range("E21").Value = "Aggiunta funzionalità cancellazione avviso"    # idx 42
$ws.Range("E22").Value = "Aggiunta affordance"                           # idx 43
$ws.Range("D34").Value = "Aggiunta affordance a tutto il menu"           # idx 44

$ws.Range("B3").Value = "Disegno xml pagina dispensa"   # idx 45
$ws.Range("B2").Value = "Dispensa"                      # idx 46

$ws.Range("D31").Value = "Aggiungi Piatto(AP)"                             # idx 47
$ws.Range("D32").Value = "AP: Inserimento dati e richieste"                # idx 48
$ws.Range("D33").Value = "AP: Inserimento automompletamento con opendata"  # idx 49

# Fill in the remaining re-used (already-existing) strings for the rest of
# the grid (rows 19-33), matching the commit's expanded todo list.
$ws.Range("D19").Value = "Struttura FXML statico"
$ws.Range("D24").Value = "Struttura FXML statico"
$ws.Range("D29").Value = "Struttura FXML statico"

$ws.Range("E20").Value = "Aggiunta card visuale dinamica"

# ---------------------------------------------------------------------------
# 2) Second table's header row (row 17) relabelled.
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = "Login/SignUp"
$ws.Range("C17").Value = "GestioneProfilo"
$ws.Range("D17").Value = "Menu"
$ws.Range("E17").Value = "Avvisi"
$ws.Range("F17").Value = "Risorse Umane"
$ws.Range("G17").Value = "PersonalizzazioneRistorante"

# ---------------------------------------------------------------------------
# 3) Clear out the old "RESPONSABILE: / Alfredo / Cristian" footer row (28)
#    and give it the "Modifica Piatto (MP)" label under column D instead.
# ---------------------------------------------------------------------------
$ws.Range("A28").ClearContents()
$ws.Range("B28").ClearContents()
$ws.Range("C28").ClearContents()
$ws.Range("E28").ClearContents()
$ws.Range("F28").ClearContents()
$ws.Range("G28").ClearContents()

# ---------------------------------------------------------------------------
# 4) Resize Tabella13 (table2) from A16:I28 to A16:I33 to cover the 5 new
#    rows, then format the new rows the same way as the old footer row.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(2)
$lo.Resize($ws.Range("A16:I33"))

$ws.Range("B29:C29").Value = ""
$ws.Range("E29:I29").Value = ""
$ws.Range("B30:C30").Value = ""
$ws.Range("E30:I30").Value = ""
$ws.Range("B31:C31").Value = ""
$ws.Range("E31:I31").Value = ""
$ws.Range("B32:C32").Value = ""
$ws.Range("E32:I32").Value = ""
$ws.Range("B33:C33").Value = ""
$ws.Range("E33:I33").Value = ""

# ---------------------------------------------------------------------------
# 5) Styling touch-ups.
# ---------------------------------------------------------------------------
# 5a. Remove the old thin top-border separator that used to sit above rows
#     13 and 28 (style now renders with no border at all).
$noLine = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$ws.Range("A13:I13").Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).LineStyle = $noLine
$ws.Range("A28:I28").Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).LineStyle = $noLine

# 5b. Give column A of the five new rows (29-33) a boxed look (medium left
#     + medium right + thin top), mirroring the new border style added to
#     styles.xml.
$continuous = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$medium = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlMedium
$thin = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
$left = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft
$right = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight
$top = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop
$bottom = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom

$boxRange = $ws.Range("A29:A33")
$boxRange.Borders.Item($left).LineStyle = $continuous
$boxRange.Borders.Item($left).Weight = $medium
$boxRange.Borders.Item($right).LineStyle = $continuous
$boxRange.Borders.Item($right).Weight = $medium
$boxRange.Borders.Item($top).LineStyle = $continuous
$boxRange.Borders.Item($top).Weight = $thin

# 5c. Underline the final "Aggiunta affordance a tutto il menu" cell (D34)
#     with a thin bottom border.
$ws.Range("D34").Borders.Item($bottom).LineStyle = $continuous
$ws.Range("D34").Borders.Item($bottom).Weight = $thin

# ---------------------------------------------------------------------------
# 6) Column D is widened to fit the new, longer descriptions.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 56.14

# ---------------------------------------------------------------------------
# 7) Move the active selection from E18 to B4 (matches the saved cursor
#    position in the target file).
# ---------------------------------------------------------------------------
$ws.Range("B4").Select() | Out-Null
